$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert "set-collaboration-enquiry-status.js" right after
#    "set-contact-name.js" (first list, numId=4), including the
#    __DdeLink__46_1190187160 bookmark that wraps the new paragraph text.
# ---------------------------------------------------------------------------
$p1 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd() -eq "set-contact-name.js") {
        $p1 = $p
        break
    }
}
$p1.Range.InsertParagraphAfter()
$newPara1 = $p1.Next()
$newPara1.Range.Text = "set-collaboration-enquiry-status.js"
$newPara1 = $p1.Next()
$zeroRange1 = $d.Range($newPara1.Range.Start, $newPara1.Range.Start)
$d.Bookmarks.Add("__DdeLink__46_1190187160", $zeroRange1)

# ---------------------------------------------------------------------------
# 2) Insert "set-collaboration-enquiry-status.js" right after
#    "Run set Searchable contact.js" (numId=3, ilvl=1), no bookmark.
# ---------------------------------------------------------------------------
$p2 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd() -eq "Run set Searchable contact.js") {
        $p2 = $p
        break
    }
}
$p2.Range.InsertParagraphAfter()
$newPara2 = $p2.Next()
$newPara2.Range.Text = "set-collaboration-enquiry-status.js"

# ---------------------------------------------------------------------------
# 3) Move first picture horizontally: posOffset 0 -> 719455 EMU
# ---------------------------------------------------------------------------
$shape1 = $d.Shapes.Item(1)
$shape1.Left = 719455 / 12700.0

# ---------------------------------------------------------------------------
# 4) Move second picture vertically: posOffset -59690 -> 59690 EMU
# ---------------------------------------------------------------------------
$shape2 = $d.Shapes.Item(2)
$shape2.Top = 59690 / 12700.0
